$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.842714
$ws.Range("H2").Value = 41.685428
$ws.Range("I2").Value = 0.06176264451223276
$ws.Range("J2").Value = 0.04208443214243528
$ws.Range("M2").Value = 181.556244
$ws.Range("N2").Value = 544.668732
$ws.Range("O2").Value = 0.393453292404907
$ws.Range("P2").Value = 0.3935455037432071
$ws.Range("Q2").Value = 3784.124868606216
$ws.Range("R2").Value = 22704.7492116373
$ws.Range("S2").Value = 0.02430071583097184
$ws.Range("T2").Value = 0.01656213904724151

$ws.Range("G3").Value = 20.842714
$ws.Range("H3").Value = 41.685428
$ws.Range("I3").Value = 0.06176264451223276
$ws.Range("J3").Value = 0.04208443214243528
$ws.Range("M3").Value = 0.324361
$ws.Range("N3").Value = 0.648722
$ws.Range("O3").Value = 0.0007029276469155644
$ws.Range("P3").Value = 0.0004687282586276696
$ws.Range("Q3").Value = 6.760563555754
$ws.Range("R3").Value = 27.042254223016
$ws.Range("S3").Value = [double]"4.341467037426627E-05"
$ws.Range("T3").Value = [double]"1.972616259345802E-05"

$ws.Range("G4").Value = 20.842714
$ws.Range("H4").Value = 41.685428
$ws.Range("I4").Value = 0.06176264451223276
$ws.Range("J4").Value = 0.04208443214243528
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 155.929759
$ws.Range("N4").Value = 467.789277
$ws.Range("O4").Value = 0.3379177477501335
$ws.Range("P4").Value = 0.3379969435488647
$ws.Range("Q4").Value = 3249.999370925927
$ws.Range("R4").Value = 19499.99622555556
$ws.Range("S4").Value = 0.02087069372866584
$ws.Range("T4").Value = 0.01422440943513272

$ws.Range("G5").Value = 20.842714
$ws.Range("H5").Value = 41.685428
$ws.Range("I5").Value = 0.06176264451223276
$ws.Range("J5").Value = 0.04208443214243528
$ws.Range("M5").Value = 123.632576
$ws.Range("N5").Value = 370.897728
$ws.Range("O5").Value = 0.2679260321980438
$ws.Range("P5").Value = 0.2679888244493004
$ws.Range("Q5").Value = 2576.838422651264
$ws.Range("R5").Value = 15461.03053590759
$ws.Range("S5").Value = 0.01654782028222081
$ws.Range("T5").Value = 0.01127815749746758

$ws.Range("I6").Value = 0.06445182531459281
$ws.Range("J6").Value = 0.0658752185158826
$ws.Range("M6").Value = 181.556244
$ws.Range("N6").Value = 544.668732
$ws.Range("O6").Value = 0.393453292404907
$ws.Range("P6").Value = 0.3935455037432071
$ws.Range("Q6").Value = 3948.887825742444
$ws.Range("R6").Value = 35539.990431682
$ws.Range("S6").Value = 0.02535878287153247
$ws.Range("T6").Value = 0.02592489605502686

$ws.Range("I7").Value = 0.06445182531459281
$ws.Range("J7").Value = 0.0658752185158826
$ws.Range("M7").Value = 0.324361
$ws.Range("N7").Value = 0.648722
$ws.Range("O7").Value = 0.0007029276469155644
$ws.Range("P7").Value = 0.0004687282586276696
$ws.Range("Q7").Value = 7.054922352577667
$ws.Range("R7").Value = 42.329534115466
$ws.Range("S7").Value = [double]"4.530496990779974E-05"
$ws.Range("T7").Value = [double]"3.087757646166687E-05"

$ws.Range("I8").Value = 0.06445182531459281
$ws.Range("J8").Value = 0.0658752185158826
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 155.929759
$ws.Range("N8").Value = 467.789277
$ws.Range("O8").Value = 0.3379177477501335
$ws.Range("P8").Value = 0.3379969435488647
$ws.Range("Q8").Value = 3391.506198960876
$ws.Range("R8").Value = 30523.55579064788
$ws.Range("S8").Value = 0.02177941564869224
$ws.Range("T8").Value = 0.02226562251398189

$ws.Range("I9").Value = 0.06445182531459281
$ws.Range("J9").Value = 0.0658752185158826
$ws.Range("M9").Value = 123.632576
$ws.Range("N9").Value = 370.897728
$ws.Range("O9").Value = 0.2679260321980438
$ws.Range("P9").Value = 0.2679888244493004
$ws.Range("Q9").Value = 2689.03543869071
$ws.Range("R9").Value = 24201.31894821639
$ws.Range("S9").Value = 0.01726832182446029
$ws.Range("T9").Value = 0.01765382237041217

$ws.Range("G10").Value = 38.15032833333333
$ws.Range("H10").Value = 114.450985
$ws.Range("I10").Value = 0.1130498248393481
$ws.Range("J10").Value = 0.115546485737591
$ws.Range("M10").Value = 181.556244
$ws.Range("N10").Value = 544.668732
$ws.Range("O10").Value = 0.393453292404907
$ws.Range("P10").Value = 0.3935455037432071
$ws.Range("Q10").Value = 6926.43031956678
$ws.Range("R10").Value = 62337.87287610102
$ws.Range("S10").Value = 0.04447982578883954
$ws.Range("T10").Value = 0.04547279993535756

$ws.Range("G11").Value = 38.15032833333333
$ws.Range("H11").Value = 114.450985
$ws.Range("I11").Value = 0.1130498248393481
$ws.Range("J11").Value = 0.115546485737591
$ws.Range("M11").Value = 0.324361
$ws.Range("N11").Value = 0.648722
$ws.Range("O11").Value = 0.0007029276469155644
$ws.Range("P11").Value = 0.0004687282586276696
$ws.Range("Q11").Value = 12.37447864852833
$ws.Range("R11").Value = 74.24687189117
$ws.Range("S11").Value = [double]"7.946584735853966E-05"
$ws.Range("T11").Value = [double]"5.415990305032791E-05"

$ws.Range("G12").Value = 38.15032833333333
$ws.Range("H12").Value = 114.450985
$ws.Range("I12").Value = 0.1130498248393481
$ws.Range("J12").Value = 0.115546485737591
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 155.929759
$ws.Range("N12").Value = 467.789277
$ws.Range("O12").Value = 0.3379177477501335
$ws.Range("P12").Value = 0.3379969435488647
$ws.Range("Q12").Value = 5948.771502787539
$ws.Range("R12").Value = 53538.94352508785
$ws.Range("S12").Value = 0.0382015421932596
$ws.Range("T12").Value = 0.03905435901711826

$ws.Range("G13").Value = 38.15032833333333
$ws.Range("H13").Value = 114.450985
$ws.Range("I13").Value = 0.1130498248393481
$ws.Range("J13").Value = 0.115546485737591
$ws.Range("M13").Value = 123.632576
$ws.Range("N13").Value = 370.897728
$ws.Range("O13").Value = 0.2679260321980438
$ws.Range("P13").Value = 0.2679888244493004
$ws.Range("Q13").Value = 4716.623367095787
$ws.Range("R13").Value = 42449.61030386209
$ws.Range("S13").Value = 0.03028899100989038
$ws.Range("T13").Value = 0.03096516688206488

$ws.Range("G14").Value = 1.0325075
$ws.Range("H14").Value = 2.065015
$ws.Range("I14").Value = 0.003059601243807028
$ws.Range("J14").Value = 0.00208478088891425
$ws.Range("M14").Value = 181.556244
$ws.Range("N14").Value = 544.668732
$ws.Range("O14").Value = 0.393453292404907
$ws.Range("P14").Value = 0.3935455037432071
$ws.Range("Q14").Value = 187.45818360183
$ws.Range("R14").Value = 1124.74910161098
$ws.Range("S14").Value = 0.001203810182822024
$ws.Range("T14").Value = 0.0008204561451219698

$ws.Range("G15").Value = 1.0325075
$ws.Range("H15").Value = 2.065015
$ws.Range("I15").Value = 0.003059601243807028
$ws.Range("J15").Value = 0.00208478088891425
$ws.Range("M15").Value = 0.324361
$ws.Range("N15").Value = 0.648722
$ws.Range("O15").Value = 0.0007029276469155644
$ws.Range("P15").Value = 0.0004687282586276696
$ws.Range("Q15").Value = 0.3349051652075
$ws.Range("R15").Value = 1.33962066083
$ws.Range("S15").Value = [double]"2.150678302809209E-06"
$ws.Range("T15").Value = [double]"9.771957156810216E-07"

$ws.Range("G16").Value = 1.0325075
$ws.Range("H16").Value = 2.065015
$ws.Range("I16").Value = 0.003059601243807028
$ws.Range("J16").Value = 0.00208478088891425
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 155.929759
$ws.Range("N16").Value = 467.789277
$ws.Range("O16").Value = 0.3379177477501335
$ws.Range("P16").Value = 0.3379969435488647
$ws.Range("Q16").Value = 160.9986456406925
$ws.Range("R16").Value = 965.991873844155
$ws.Range("S16").Value = 0.001033893561320778
$ws.Range("T16").Value = 0.0007046495684221018

$ws.Range("G17").Value = 1.0325075
$ws.Range("H17").Value = 2.065015
$ws.Range("I17").Value = 0.003059601243807028
$ws.Range("J17").Value = 0.00208478088891425
$ws.Range("M17").Value = 123.632576
$ws.Range("N17").Value = 370.897728
$ws.Range("O17").Value = 0.2679260321980438
$ws.Range("P17").Value = 0.2679888244493004
$ws.Range("Q17").Value = 127.65156196432
$ws.Range("R17").Value = 765.90937178592
$ws.Range("S17").Value = 0.0008197468213614167
$ws.Range("T17").Value = 0.0005586979796544974

$ws.Range("G18").Value = 239.6229553333334
$ws.Range("H18").Value = 718.868866
$ws.Range("I18").Value = 0.710068151739898
$ws.Range("J18").Value = 0.7257497274703861
$ws.Range("M18").Value = 181.556244
$ws.Range("N18").Value = 544.668732
$ws.Range("O18").Value = 0.393453292404907
$ws.Range("P18").Value = 0.3935455037432071
$ws.Range("Q18").Value = 43505.04374649977
$ws.Range("R18").Value = 391545.3937184979
$ws.Range("S18").Value = 0.27937865213393
$ws.Range("T18").Value = 0.2856155420888284

$ws.Range("G19").Value = 239.6229553333334
$ws.Range("H19").Value = 718.868866
$ws.Range("I19").Value = 0.710068151739898
$ws.Range("J19").Value = 0.7257497274703861
$ws.Range("M19").Value = 0.324361
$ws.Range("N19").Value = 0.648722
$ws.Range("O19").Value = 0.0007029276469155644
$ws.Range("P19").Value = 0.0004687282586276696
$ws.Range("Q19").Value = 77.72434141487534
$ws.Range("R19").Value = 466.346048489252
$ws.Range("S19").Value = 0.0004991265350522105
$ws.Range("T19").Value = 0.0003401794059566999

$ws.Range("G20").Value = 239.6229553333334
$ws.Range("H20").Value = 718.868866
$ws.Range("I20").Value = 0.710068151739898
$ws.Range("J20").Value = 0.7257497274703861
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 155.929759
$ws.Range("N20").Value = 467.789277
$ws.Range("O20").Value = 0.3379177477501335
$ws.Range("P20").Value = 0.3379969435488647
$ws.Range("Q20").Value = 37364.34967599444
$ws.Range("R20").Value = 336279.1470839499
$ws.Range("S20").Value = 0.2399446305850464
$ws.Range("T20").Value = 0.245301189666412

$ws.Range("G21").Value = 239.6229553333334
$ws.Range("H21").Value = 718.868866
$ws.Range("I21").Value = 0.710068151739898
$ws.Range("J21").Value = 0.7257497274703861
$ws.Range("M21").Value = 123.632576
$ws.Range("N21").Value = 370.897728
$ws.Range("O21").Value = 0.2679260321980438
$ws.Range("P21").Value = 0.2679888244493004
$ws.Range("Q21").Value = 29625.20323659294
$ws.Range("R21").Value = 266626.8291293365
$ws.Range("S21").Value = 0.1902457424858694
$ws.Range("T21").Value = 0.1944928163091889

$ws.Range("G22").Value = 16.06600466666667
$ws.Range("H22").Value = 48.198014
$ws.Range("I22").Value = 0.04760795235012129
$ws.Range("J22").Value = 0.04865935524479072
$ws.Range("M22").Value = 181.556244
$ws.Range("N22").Value = 544.668732
$ws.Range("O22").Value = 0.393453292404907
$ws.Range("P22").Value = 0.3935455037432071
$ws.Range("Q22").Value = 2916.883463366472
$ws.Range("R22").Value = 26251.95117029825
$ws.Range("S22").Value = 0.01873150559681115
$ws.Range("T22").Value = 0.01914967047163083

$ws.Range("G23").Value = 16.06600466666667
$ws.Range("H23").Value = 48.198014
$ws.Range("I23").Value = 0.04760795235012129
$ws.Range("J23").Value = 0.04865935524479072
$ws.Range("M23").Value = 0.324361
$ws.Range("N23").Value = 0.648722
$ws.Range("O23").Value = 0.0007029276469155644
$ws.Range("P23").Value = 0.0004687282586276696
$ws.Range("Q23").Value = 5.211185339684667
$ws.Range("R23").Value = 31.267112038108
$ws.Range("S23").Value = [double]"3.346494591993907E-05"
$ws.Range("T23").Value = [double]"2.280801484983591E-05"

$ws.Range("G24").Value = 16.06600466666667
$ws.Range("H24").Value = 48.198014
$ws.Range("I24").Value = 0.04760795235012129
$ws.Range("J24").Value = 0.04865935524479072
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 155.929759
$ws.Range("N24").Value = 467.789277
$ws.Range("O24").Value = 0.3379177477501335
$ws.Range("P24").Value = 0.3379969435488647
$ws.Range("Q24").Value = 2505.168235766209
$ws.Range("R24").Value = 22546.51412189588
$ws.Range("S24").Value = 0.01608757203314866
$ws.Range("T24").Value = 0.01644671334779768

$ws.Range("G25").Value = 16.06600466666667
$ws.Range("H25").Value = 48.198014
$ws.Range("I25").Value = 0.04760795235012129
$ws.Range("J25").Value = 0.04865935524479072
$ws.Range("M25").Value = 123.632576
$ws.Range("N25").Value = 370.897728
$ws.Range("O25").Value = 0.2679260321980438
$ws.Range("P25").Value = 0.2679888244493004
$ws.Range("Q25").Value = 1986.281542968022
$ws.Range("R25").Value = 17876.53388671219
$ws.Range("S25").Value = 0.01275540977424153
$ws.Range("T25").Value = 0.01304016341051236
